$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 27 ("POU hot water heaters") and shift the rows below it up
$ws.Rows("27").Delete()

# Update the selection to reflect the new active cell after edit
$ws.Range("C33").Select()
